$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the "D3 D5" / 1N5711 Schottky-diode row (row 7). ---
# Deleting the row shifts every row below it up by one, which matches the
# rest of the diff (dimension A1:U25 -> A1:U24, shared formula range
# E3:E25 -> E3:E24, and every row 9..25 sliding into 8..24).
$ws.Rows(7).Delete()

# --- 2. The row that used to be row 8 ("D4 D6" / 1N4148) is now row 7. ---
# It absorbs the deleted part's reference designators and quantity, so it
# becomes "D3 D4 D5 D6" with Qty/Have = 4 (was 2).
$ws.Range("A7").Value = "D3 D4 D5 D6"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 4

# --- 3. Conditional formatting cleanup. ---
# There were two cellIs > 0 rules: one on E2:E25 (red "needs more" highlight)
# and a stray leftover one on H7 (same red style, no real cell content).
# Drop the big rule and re-target the H7 rule (reusing its dxf) to the new
# E2:E24 range, so only one rule/dxf stays in active use.
$fcBig = $ws.Range("E2:E25").FormatConditions.Item(1)
$fcH7 = $ws.Range("H7").FormatConditions.Item(1)
$fcBig.Delete()
$fcH7.ModifyAppliesToRange($ws.Range("E2:E24"))

# --- 4. Move the "D-shaft for small knobs" comment from B23 to B22. ---
# Row content shifted up with the delete, but comment anchors don't follow
# automatically, so re-create it one row up with the same text.
$oldText = $ws.Range("B23").Comment.Text()
$ws.Range("B23").Comment.Delete()
$ws.Range("B22").AddComment($oldText) | Out-Null
